$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain numeric-looking string must be forced to
# stay as TEXT (matching the source file, which stores every data cell as an
# inline string) -- otherwise the COM layer auto-converts them to numbers.

$ws.Range("D2").Value = "59.808.96"
$ws.Range("E2").Value = "  +0.90%  "
$ws.Range("D3").Value = "2.664.98"
$ws.Range("E3").Value = "  +2.62%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "538.33"
$ws.Range("E5").Value = "  +0.52%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.68"
$ws.Range("E6").Value = "  +3.36%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  +0.18%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.573"
$ws.Range("E8").Value = "  +0.94%  "
$ws.Range("D9").Value = "2.666.05"
$ws.Range("E9").Value = "  +2.06%  "
$ws.Range("E10").Value = "  +3.26%  "
$ws.Range("E11").Value = "  +0.75%  "
$ws.Range("E12").Value = "  +1.09%  "
$ws.Range("E13").Value = "  -0.75%  "
$ws.Range("D14").Value = "3.131.93"
$ws.Range("E14").Value = "  +2.44%  "
$ws.Range("D15").Value = "59.726.95"
$ws.Range("E15").Value = "  +0.88%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "21.21"
$ws.Range("E16").Value = "  +2.92%  "
$ws.Range("D17").Value = "2.677.62"
$ws.Range("E17").Value = "  +1.02%  "
$ws.Range("E18").Value = "  +1.20%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "345.49"
$ws.Range("E19").Value = "  -0.23%  "
$ws.Range("E20").Value = "  +1.92%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.41"
$ws.Range("E21").Value = "  +2.92%  "
$ws.Range("E22").Value = "  -0.26%  "
$ws.Range("E23").Value = "  +0.20%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.60"
$ws.Range("E24").Value = "  -0.75%  "
$ws.Range("E25").Value = "  +2.14%  "
$ws.Range("E26").Value = "  -1.51%  "
$ws.Range("E27").Value = "  -0.04%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.31"
$ws.Range("E28").Value = "  +1.16%  "
$ws.Range("D29").Value = "0.0₃0754"
$ws.Range("E29").Value = "  +1.90%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.999"
$ws.Range("E30").Value = "  +0.08%  "
$ws.Range("E31").Value = "  +1.68%  "
$ws.Range("E32").Value = "  +0.48%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.05"
$ws.Range("E33").Value = "  +1.16%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "150.37"
$ws.Range("E34").Value = "  +0.53%  "
$ws.Range("E35").Value = "  +1.06%  "
$ws.Range("E36").Value = "  +1.95%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.845"
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("E38").Value = "  -1.05%  "
$ws.Range("E39").Value = "  -0.83%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "293.14"
$ws.Range("E40").Value = "  +5.04%  "
$ws.Range("E41").Value = "  +1.95%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.24%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.605"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "19.50"
$ws.Range("E44").Value = "  +5.17%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0542"
$ws.Range("E45").Value = "  +4.07%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "10.73"
$ws.Range("E47").Value = "  -1.30%  "
$ws.Range("D48").Value = "1.978.25"
$ws.Range("E48").Value = "  +1.83%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0227"
$ws.Range("E49").Value = "  +1.77%  "
$ws.Range("B50").Value = "InjectiveProtocol"
$ws.Range("C50").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "18.40"
$ws.Range("E50").Value = "  +0.10%  "
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "4.54"
$ws.Range("E51").Value = "  +0.75%  "
